$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('A2').Value = 'Fulham v Arsenal'
$ws.Range('B2').Value = 'Arsenal'
$ws.Range('C2').Value = 'England Premier League'
$ws.Range('D2').Value = '2025-10-18T16:30:00.000Z'
$ws.Range('E2').Value = '88/113 Win Tips'
$ws.Range('F2').Value = "'78"
$ws.Range('F2').Style = 'Normal'
$ws.Range('G2').Value = "'1.50"
$ws.Range('G2').Style = 'Normal'

# Row 3
$ws.Range('A3').Value = 'Torino v Napoli'
$ws.Range('B3').Value = 'Napoli'
$ws.Range('C3').Value = 'Italy Serie A'
$ws.Range('D3').Value = '2025-10-18T16:00:00.000Z'
$ws.Range('E3').Value = '33/37 Win Tips'
$ws.Range('F3').Value = "'89"
$ws.Range('F3').Style = 'Normal'
$ws.Range('G3').Value = "'1.83"
$ws.Range('G3').Style = 'Normal'

# Row 4
$ws.Range('A4').Value = 'Marseille v Le Havre'
$ws.Range('B4').Value = 'Marseille'
$ws.Range('C4').Value = 'France Ligue 1'
$ws.Range('D4').Value = '2025-10-18T19:05:00.000Z'
$ws.Range('E4').Value = '30/31 Win Tips'
$ws.Range('F4').Value = "'97"
$ws.Range('F4').Style = 'Normal'
$ws.Range('G4').Value = "'1.38"
$ws.Range('G4').Style = 'Normal'

# Row 5
$ws.Range('A5').Value = 'Atletico Madrid v Osasuna'
$ws.Range('B5').Value = 'Atletico Madrid'
$ws.Range('C5').Value = 'Spain Primera Liga'
$ws.Range('D5').Value = '2025-10-18T19:00:00.000Z'
$ws.Range('E5').Value = '27/30 Win Tips'
$ws.Range('F5').Value = "'90"
$ws.Range('F5').Style = 'Normal'
$ws.Range('G5').Value = "'1.38"
$ws.Range('G5').Style = 'Normal'

# Row 6
$ws.Range('A6').Value = 'Liverpool v Man Utd'
$ws.Range('B6').Value = 'Liverpool'
$ws.Range('C6').Value = 'England Premier League'
$ws.Range('D6').Value = '2025-10-19T15:30:00.000Z'
$ws.Range('E6').Value = '27/38 Win Tips'
$ws.Range('F6').Value = "'71"
$ws.Range('F6').Style = 'Normal'
$ws.Range('G6').Value = "'1.62"
$ws.Range('G6').Style = 'Normal'

# Row 7
$ws.Range('A7').Value = 'Leicester v Portsmouth'
$ws.Range('B7').Value = 'Leicester'
$ws.Range('C7').Value = 'England Championship'
$ws.Range('D7').Value = '2025-10-18T18:45:00.000Z'
$ws.Range('E7').Value = '22/26 Win Tips'
$ws.Range('F7').Value = "'85"
$ws.Range('F7').Style = 'Normal'
$ws.Range('G7').Value = "'1.62"
$ws.Range('G7').Style = 'Normal'

# Row 8
$ws.Range('A8').Value = 'Bayern Munich v Borussia Dortmund'
$ws.Range('B8').Value = 'Bayern Munich'
$ws.Range('C8').Value = 'Germany Bundesliga I'
$ws.Range('D8').Value = '2025-10-18T16:30:00.000Z'
$ws.Range('E8').Value = '16/22 Win Tips'
$ws.Range('F8').Value = "'73"
$ws.Range('F8').Style = 'Normal'
$ws.Range('G8').Value = "'1.40"
$ws.Range('G8').Style = 'Normal'

# Row 9
$ws.Range('A9').Value = 'Tottenham v Aston Villa'
$ws.Range('B9').Value = 'Tottenham'
$ws.Range('C9').Value = 'England Premier League'
$ws.Range('D9').Value = '2025-10-19T13:00:00.000Z'
$ws.Range('E9').Value = '14/26 Win Tips'
$ws.Range('F9').Value = "'54"
$ws.Range('F9').Style = 'Normal'
$ws.Range('G9').Value = "'2.15"
$ws.Range('G9').Style = 'Normal'

# Row 10
$ws.Range('A10').Value = 'Angers v Monaco'
$ws.Range('B10').Value = 'Monaco'
$ws.Range('C10').Value = 'France Ligue 1'
$ws.Range('D10').Value = '2025-10-18T17:00:00.000Z'
$ws.Range('E10').Value = '13/17 Win Tips'
$ws.Range('F10').Value = "'76"
$ws.Range('F10').Style = 'Normal'
$ws.Range('G10').Value = "'1.60"
$ws.Range('G10').Style = 'Normal'

# Row 11
$ws.Range('A11').Value = 'Dinamo Zagreb v NK Osijek'
$ws.Range('B11').Value = 'Dinamo Zagreb'
$ws.Range('C11').Value = 'Croatia HNL'
$ws.Range('D11').Value = '2025-10-18T16:00:00.000Z'
$ws.Range('E11').Value = '12/16 Win Tips'
$ws.Range('F11').Value = "'75"
$ws.Range('F11').Style = 'Normal'
$ws.Range('G11').Value = "'1.36"
$ws.Range('G11').Style = 'Normal'

# Row 12
$ws.Range('A12').Value = 'OH Leuven v Club Brugge'
$ws.Range('B12').Value = 'Club Brugge'
$ws.Range('C12').Value = 'Belgium First Division A'
$ws.Range('D12').Value = '2025-10-18T16:15:00.000Z'
$ws.Range('E12').Value = '12/14 Win Tips'
$ws.Range('F12').Value = "'86"
$ws.Range('F12').Style = 'Normal'
$ws.Range('G12').Value = "'1.67"
$ws.Range('G12').Style = 'Normal'

# Row 13
$ws.Range('A13').Value = 'Getafe v Real Madrid'
$ws.Range('B13').Value = 'Real Madrid'
$ws.Range('C13').Value = 'Spain Primera Liga'
$ws.Range('D13').Value = '2025-10-19T19:00:00.000Z'
$ws.Range('E13').Value = '12/12 Win Tips'
$ws.Range('F13').Value = "'100"
$ws.Range('F13').Style = 'Normal'
$ws.Range('G13').Value = "'1.45"
$ws.Range('G13').Style = 'Normal'

# Row 14
$ws.Range('A14').Value = 'Villarreal v Real Betis'
$ws.Range('B14').Value = 'Yes'
$ws.Range('C14').Value = 'Spain Primera Liga'
$ws.Range('D14').Value = '2025-10-18T16:30:00.000Z'
$ws.Range('E14').Value = '10/10 Win Tips'
$ws.Range('F14').Value = "'100"
$ws.Range('F14').Style = 'Normal'
$ws.Range('G14').Value = "'1.57"
$ws.Range('G14').Style = 'Normal'

# Row 15
$ws.Range('A15').Value = 'PSV v Go Ahead Eagles'
$ws.Range('B15').Value = 'PSV'
$ws.Range('C15').Value = 'Netherlands Eredivisie'
$ws.Range('D15').Value = '2025-10-18T16:45:00.000Z'
$ws.Range('E15').Value = '10/11 Win Tips'
$ws.Range('F15').Value = "'91"
$ws.Range('F15').Style = 'Normal'
$ws.Range('G15').Value = "'1.25"
$ws.Range('G15').Style = 'Normal'

# Row 16
$ws.Range('A16').Value = 'Union Saint Gilloise v Charleroi'
$ws.Range('B16').Value = 'Union Saint Gilloise'
$ws.Range('C16').Value = 'Belgium First Division A'
$ws.Range('D16').Value = '2025-10-18T18:45:00.000Z'
$ws.Range('E16').Value = '10/10 Win Tips'
$ws.Range('F16').Value = "'100"
$ws.Range('F16').Style = 'Normal'
$ws.Range('G16').Value = "'1.45"
$ws.Range('G16').Style = 'Normal'

# Row 17
$ws.Range('A17').Value = 'SK Brann v Haugesund'
$ws.Range('B17').Value = 'SK Brann'
$ws.Range('C17').Value = 'Norway Eliteserien'
$ws.Range('D17').Value = '2025-10-18T16:00:00.000Z'
$ws.Range('E17').Value = '9/9 Win Tips'
$ws.Range('F17').Value = "'100"
$ws.Range('F17').Style = 'Normal'
$ws.Range('G17').Value = "'1.11"
$ws.Range('G17').Style = 'Normal'

# Row 18
$ws.Range('A18').Value = 'Istanbul Basaksehir v Galatasaray'
$ws.Range('B18').Value = 'Galatasaray'
$ws.Range('C18').Value = 'Turkey Super Lig'
$ws.Range('D18').Value = '2025-10-18T17:00:00.000Z'
$ws.Range('E18').Value = '9/11 Win Tips'
$ws.Range('F18').Value = "'82"
$ws.Range('F18').Style = 'Normal'
$ws.Range('G18').Value = "'1.70"
$ws.Range('G18').Style = 'Normal'

# Row 19
$ws.Range('A19').Value = 'Basel v Winterthur'
$ws.Range('B19').Value = 'Basel'
$ws.Range('C19').Value = 'Switzerland Super League'
$ws.Range('D19').Value = '2025-10-18T16:00:00.000Z'
$ws.Range('E19').Value = '7/8 Win Tips'
$ws.Range('F19').Value = "'88"
$ws.Range('F19').Style = 'Normal'
$ws.Range('G19').Value = "'1.27"
$ws.Range('G19').Style = 'Normal'

# Row 20
$ws.Range('A20').Value = 'Slavia Prague v FC Zlin'
$ws.Range('B20').Value = 'Slavia Prague'
$ws.Range('C20').Value = 'Czech Republic First League'
$ws.Range('D20').Value = '2025-10-18T16:00:00.000Z'
$ws.Range('E20').Value = '6/8 Win Tips'
$ws.Range('F20').Value = "'75"
$ws.Range('F20').Style = 'Normal'
$ws.Range('G20').Value = "'1.18"
$ws.Range('G20').Style = 'Normal'

# Row 21
$ws.Range('A21').Value = 'Roma v Inter Milan'
$ws.Range('B21').Value = 'Draw'
$ws.Range('C21').Value = 'Italy Serie A'
$ws.Range('D21').Value = '2025-10-18T18:45:00.000Z'
$ws.Range('E21').Value = '6/14 Win Tips'
$ws.Range('F21').Value = "'43"
$ws.Range('F21').Style = 'Normal'
$ws.Range('G21').Value = "'3.30"
$ws.Range('G21').Style = 'Normal'
